$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86-164 down to 87-165.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new data record.
$ws.Cells.Item(86, 1).Value = 8
$ws.Cells.Item(86, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(86, 3).Value = "Coquimbo"
$ws.Cells.Item(86, 4).Value = $ws.Cells.Item(87, 4).Value()
$ws.Cells.Item(86, 4).NumberFormat = $ws.Cells.Item(87, 4).NumberFormat()
$ws.Cells.Item(86, 5).Value = 4
$ws.Cells.Item(86, 6).Value = "Fruta"
$ws.Cells.Item(86, 7).Value = 100109
$ws.Cells.Item(86, 8).Value = "Uva"
$ws.Cells.Item(86, 9).Value = 100109001
$ws.Cells.Item(86, 10).Value = "Uva"
$ws.Cells.Item(86, 11).Value = "Moscatel rosada"
$ws.Cells.Item(86, 12).Value = "Primera"
$ws.Cells.Item(86, 13).Value = 400
$ws.Cells.Item(86, 14).Value = 11000
$ws.Cells.Item(86, 15).Value = 12000
$ws.Cells.Item(86, 16).Value = 11500
$ws.Cells.Item(86, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(86, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(86, 19).Value = 1150
$ws.Cells.Item(86, 20).Value = 10

Write-Host "done"
